$wb = $excel.ActiveWorkbook

# --- SNAGS sheet: insert a new column before column A, shifting all data
#     from A:B to B:C, and update the selection on that sheet ---
$snags = $wb.Worksheets.Item("SNAGS")
$snags.Columns.Item(1).Insert()
$snags.Range("C6").Select()

# --- IO sheet: becomes the active/selected sheet, with a new active cell ---
$io = $wb.Worksheets.Item("IO")
$io.Activate()
$io.Range("K27").Select()
